$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 12, shifting existing rows 12-14 down to 13-15.
$ws.Rows.Item(12).Insert()

# Populate the new row 12 with the new weekly data point.
# (Other columns mirror the constant values used throughout this sheet.)
$ws.Cells.Item(12, 1).Value = 7
$ws.Cells.Item(12, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(12, 3).Value = "Ñuble"
$ws.Cells.Item(12, 4).Value = 44813
$ws.Cells.Item(12, 5).Value = 16
$ws.Cells.Item(12, 6).Value = 100112012
$ws.Cells.Item(12, 7).Value = "Espinaca"
$ws.Cells.Item(12, 8).Value = "Sin especificar"
$ws.Cells.Item(12, 9).Value = "Primera"
$ws.Cells.Item(12, 10).Value = 120
$ws.Cells.Item(12, 11).Value = 7000
$ws.Cells.Item(12, 12).Value = 7500
$ws.Cells.Item(12, 13).Value = 7250
$ws.Cells.Item(12, 14).Value = "$/cuna 10 kilos"
$ws.Cells.Item(12, 15).Value = "Provincia de Diguillín"
$ws.Cells.Item(12, 16).Value = 725
$ws.Cells.Item(12, 17).Value = 10
$ws.Cells.Item(12, 18).Value = "Hortaliza"
